$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1081400015761023
$ws.Range("H2").Value = -17.71466198968813
$ws.Range("I2").Value = 16.52007898291481

$ws.Range("G3").Value = 0.1189519931576234
$ws.Range("H3").Value = 33.5980894574028

$ws.Range("G4").Value = -0.6879409537913153
$ws.Range("H4").Value = -11.79788612917265

$ws.Range("G5").Value = -0.627564809514831
$ws.Range("H5").Value = -2.779520555373815

$ws.Range("G6").Value = 0.1370603655286926
$ws.Range("H6").Value = -44.30880504169735

$ws.Range("G7").Value = 0.3488370368176556
$ws.Range("H7").Value = 112.9366988113364

$ws.Range("G8").Value = 0.1119416242393618
$ws.Range("H8").Value = -32.283376746107

$ws.Range("G9").Value = 0.227389036905903
$ws.Range("H9").Value = 16.55445061663369

$ws.Range("G10").Value = -0.1220396634835586
$ws.Range("H10").Value = -113.5563663106809

$ws.Range("G11").Value = -0.1243865623968502
$ws.Range("H11").Value = -4.730419050405798

$ws.Range("G12").Value = 0.2060282259208529
$ws.Range("H12").Value = 29.55363843344319

$ws.Range("G13").Value = 0.2251849522034897
$ws.Range("H13").Value = 9.492871260575503

$ws.Range("G14").Value = 0.1892185175612739
$ws.Range("H14").Value = -0.08161601034078958

$ws.Range("G15").Value = 0.2511973853183916
$ws.Range("H15").Value = 0.5169826697811607

$ws.Range("G16").Value = 0.03292401232602186
$ws.Range("H16").Value = -9.748384046911989

$ws.Range("G17").Value = 0.05410248177144972
$ws.Range("H17").Value = 52.52771859302774

$ws.Range("G18").Value = 0.0850246139402005
$ws.Range("H18").Value = -50.94192947167575

$ws.Range("G19").Value = 0.1779039450263613
$ws.Range("H19").Value = 41.47099827149915

$ws.Range("G20").Value = 0.09123014398003387
$ws.Range("H20").Value = -20.42865075315715

$ws.Range("G21").Value = 0.1729699088839599
$ws.Range("H21").Value = 72.28280038365934

$ws.Range("G22").Value = 0.05795863911377596
$ws.Range("H22").Value = -38.47090708946001

$ws.Range("G23").Value = 0.1551288042547354
$ws.Range("H23").Value = 42.98905115971322

$ws.Range("G24").Value = -0.1669877398190396
$ws.Range("H24").Value = -33.99287002872365

$ws.Range("G25").Value = -0.1517928825179598
$ws.Range("H25").Value = 31.76439922780055

$ws.Range("G26").Value = 0.2068856818592077
$ws.Range("H26").Value = 30.13166212237345

$ws.Range("G27").Value = 0.2358641137107203
$ws.Range("H27").Value = 17.65771646757017

$ws.Range("G28").Value = 0.008532010211720339
$ws.Range("H28").Value = 206.2457270506428

$ws.Range("G29").Value = 0.03777862236809996
$ws.Range("H29").Value = 145.6748002195506
